$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.861.33"
$ws.Range("E2").Value = "  -3.83%  "
$ws.Range("D3").Value = "3.500.50"
$ws.Range("E3").Value = "  -4.31%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'578.57"
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("D6").Value = "'175.83"
$ws.Range("E6").Value = "  -2.22%  "
$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "3.492.56"
$ws.Range("E8").Value = "  -4.31%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  -7.12%  "
$ws.Range("D11").Value = "'6.60"
$ws.Range("E11").Value = "  +6.89%  "
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("D13").Value = "'47.22"
$ws.Range("E13").Value = "  -5.19%  "
$ws.Range("E14").Value = "  -3.38%  "
$ws.Range("D15").Value = "'681.68"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "'8.87"
$ws.Range("E16").Value = "  -1.64%  "
$ws.Range("D17").Value = "4.057.21"
$ws.Range("E17").Value = "  -4.00%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.498.49"
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "68.803.08"
$ws.Range("E19").Value = "  -3.99%  "
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("D21").Value = "'17.54"
$ws.Range("E21").Value = "  -3.90%  "
$ws.Range("D22").Value = "'11.19"
$ws.Range("E22").Value = "  -3.90%  "
$ws.Range("E23").Value = "  -4.11%  "
$ws.Range("D24").Value = "'16.35"
$ws.Range("E24").Value = "  -8.53%  "
$ws.Range("D25").Value = "'98.22"
$ws.Range("E25").Value = "  -5.05%  "
$ws.Range("E26").Value = "  -4.44%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("E28").Value = "  -6.78%  "
$ws.Range("D29").Value = "'9.40"
$ws.Range("E29").Value = "  -7.89%  "
$ws.Range("D30").Value = "'32.97"
$ws.Range("E30").Value = "  -6.73%  "
$ws.Range("D31").Value = "'8.75"
$ws.Range("E31").Value = "  -4.80%  "
$ws.Range("E32").Value = "  -7.78%  "
$ws.Range("D33").Value = "'7.37"
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("E34").Value = "  -5.93%  "
$ws.Range("D35").Value = "'570.39"
$ws.Range("E35").Value = "  -1.76%  "
$ws.Range("D36").Value = "'3.63"
$ws.Range("E36").Value = "  -14.91%  "
$ws.Range("D37").Value = "'10.95"
$ws.Range("E37").Value = "  -3.50%  "
$ws.Range("E38").Value = "  -3.17%  "
$ws.Range("D39").Value = "'56.80"
$ws.Range("E39").Value = "  -5.10%  "
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").Value = "  -4.78%  "
$ws.Range("E42").Value = "  -4.40%  "
$ws.Range("D43").Value = "'0.337"
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("D44").Value = "3.418.57"
$ws.Range("E44").Value = "  -8.64%  "
$ws.Range("D45").Value = "'33.45"
$ws.Range("E45").Value = "  -6.04%  "
$ws.Range("D46").Value = "0.0₃0702"
$ws.Range("E46").Value = "  -8.49%  "
$ws.Range("D47").Value = "'2.90"
$ws.Range("E47").Value = "  +3.14%  "
$ws.Range("E48").Value = "  -7.12%  "
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").Value = "'133.97"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("E51").Value = "  -0.28%  "
